$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before FE (column 161), shifting existing FE.. to the right by 2
$ws.Range("FE1:FF1").EntireColumn.Insert()

$ws.Range("FE1").Value = "identificador_muestra"
$ws.Range("FF1").Value = "temperatura_muestra"
$ws.Range("FE1:FF1").Font.Bold = $true
